$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

$ws.Range("E2").Value = 0.001334638313017278
$ws.Range("E3").Value = 0.003491755577109679
$ws.Range("E4").Value = 0.0009435742592942553
$ws.Range("E5").Value = -0.01371428571428568
$ws.Range("E6").Value = 0.01549517179429594
$ws.Range("E7").Value = -0.002792586951002707
$ws.Range("E8").Value = -0.003382663847779965
$ws.Range("E9").Value = 0.00850749029036435
$ws.Range("E10").Value = 0.003921568627450966
$ws.Range("E11").Value = 0.003355704697986628
$ws.Range("E12").Value = 0.0088592067602562
$ws.Range("E13").Value = 0.002921966311447255
$ws.Range("E14").Value = 0.0002287805994050096
$ws.Range("E15").Value = 0.003252574138116859

$ws.Protect()
